# This script applies textual updates to the three-digit division answer
# document, replacing the old "dividend÷divisor=quotient, remainder" text
# runs with their new counterparts, per the commit diff.
#
# Each old value is unique within the document, so a simple
# Find/Replace (non-wildcard, case-sensitive, replace-all) pass per pair
# is sufficient and safe.

$d = $word.ActiveDocument

$d.Content.Find.Execute("949÷2=474, 1", $false, $false, $false, $false, $false, $true, 1, $false, "589÷5=117, 4", 2) | Out-Null
$d.Content.Find.Execute("386÷5=77, 1", $false, $false, $false, $false, $false, $true, 1, $false, "647÷7=92, 3", 2) | Out-Null
$d.Content.Find.Execute("114÷4=28, 2", $false, $false, $false, $false, $false, $true, 1, $false, "636÷6=106, 0", 2) | Out-Null
$d.Content.Find.Execute("554÷8=69, 2", $false, $false, $false, $false, $false, $true, 1, $false, "310÷8=38, 6", 2) | Out-Null
$d.Content.Find.Execute("886÷9=98, 4", $false, $false, $false, $false, $false, $true, 1, $false, "338÷9=37, 5", 2) | Out-Null
$d.Content.Find.Execute("465÷9=51, 6", $false, $false, $false, $false, $false, $true, 1, $false, "692÷6=115, 2", 2) | Out-Null
$d.Content.Find.Execute("668÷8=83, 4", $false, $false, $false, $false, $false, $true, 1, $false, "423÷7=60, 3", 2) | Out-Null
$d.Content.Find.Execute("260÷9=28, 8", $false, $false, $false, $false, $false, $true, 1, $false, "259÷7=37, 0", 2) | Out-Null
$d.Content.Find.Execute("355÷3=118, 1", $false, $false, $false, $false, $false, $true, 1, $false, "415÷4=103, 3", 2) | Out-Null
$d.Content.Find.Execute("752÷2=376, 0", $false, $false, $false, $false, $false, $true, 1, $false, "369÷9=41, 0", 2) | Out-Null
$d.Content.Find.Execute("779÷3=259, 2", $false, $false, $false, $false, $false, $true, 1, $false, "141÷8=17, 5", 2) | Out-Null
$d.Content.Find.Execute("257÷8=32, 1", $false, $false, $false, $false, $false, $true, 1, $false, "597÷7=85, 2", 2) | Out-Null
$d.Content.Find.Execute("684÷8=85, 4", $false, $false, $false, $false, $false, $true, 1, $false, "587÷5=117, 2", 2) | Out-Null
$d.Content.Find.Execute("165÷7=23, 4", $false, $false, $false, $false, $false, $true, 1, $false, "703÷9=78, 1", 2) | Out-Null
$d.Content.Find.Execute("214÷2=107, 0", $false, $false, $false, $false, $false, $true, 1, $false, "114÷3=38, 0", 2) | Out-Null
$d.Content.Find.Execute("831÷4=207, 3", $false, $false, $false, $false, $false, $true, 1, $false, "299÷9=33, 2", 2) | Out-Null
$d.Content.Find.Execute("936÷9=104, 0", $false, $false, $false, $false, $false, $true, 1, $false, "420÷9=46, 6", 2) | Out-Null
$d.Content.Find.Execute("357÷4=89, 1", $false, $false, $false, $false, $false, $true, 1, $false, "544÷9=60, 4", 2) | Out-Null
$d.Content.Find.Execute("453÷4=113, 1", $false, $false, $false, $false, $false, $true, 1, $false, "432÷6=72, 0", 2) | Out-Null
$d.Content.Find.Execute("155÷2=77, 1", $false, $false, $false, $false, $false, $true, 1, $false, "474÷6=79, 0", 2) | Out-Null
$d.Content.Find.Execute("461÷2=230, 1", $false, $false, $false, $false, $false, $true, 1, $false, "708÷7=101, 1", 2) | Out-Null
$d.Content.Find.Execute("224÷8=28, 0", $false, $false, $false, $false, $false, $true, 1, $false, "914÷9=101, 5", 2) | Out-Null
$d.Content.Find.Execute("842÷3=280, 2", $false, $false, $false, $false, $false, $true, 1, $false, "827÷2=413, 1", 2) | Out-Null
$d.Content.Find.Execute("408÷6=68, 0", $false, $false, $false, $false, $false, $true, 1, $false, "453÷3=151, 0", 2) | Out-Null
$d.Content.Find.Execute("184÷7=26, 2", $false, $false, $false, $false, $false, $true, 1, $false, "137÷9=15, 2", 2) | Out-Null
